# Fix property_category values that were left as "land" on the wrong
# sheets when the workbook was generated:
#   - 建物 (building) sheet: property_category column (I) should read
#     "building", not "land", for every data row.
#   - 汽車 (car) sheet: property_category column (H) should read "car",
#     not "land", for every data row.
$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I9").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H3").Value = "car"
